$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column at G for "mobile #" (shifts old G..L to H..M) ---
$ws.Columns("G:G").Insert()
$ws.Columns("G:G").ColumnWidth = 9

# --- Fill in row 7 (second bot) content, in the order that makes the ---
# --- shared-strings table land in the same order as the target file  ---
$ws.Range("D7").Value2 = "Meow0004"

$ws.Range("B7").Value2 = "doodlebob0045@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:doodlebob0045@gmail.com")

$ws.Range("H7").Value2 = "oedC7K7pjmlWtWYvmDQfMWgsY"
$ws.Range("H2").Copy()
$ws.Range("H7").PasteSpecial(-4122)

$ws.Range("I7").Value2 = "UpN7lmGOfZt2iXNlxaQ3FmtkZO6GjFFDnOq7Y9c5OWBCgEnchY"
$ws.Range("I2").Copy()
$ws.Range("I7").PasteSpecial(-4122)

$ws.Range("J7").Value2 = "946075554050183169-2VoOR8pk64CBXyjIU04QcnZbfSdxHzn"
$ws.Range("J2").Copy()
$ws.Range("J7").PasteSpecial(-4122)

$ws.Range("K7").Value2 = "FRouNKabtdALXLpGmSRU4raihlPVkKbm4dmoAk2L6KkDh"
$ws.Range("K2").Copy()
$ws.Range("K7").PasteSpecial(-4122)

# --- New "mobile #" header + value ---
$ws.Range("G1").Value2 = "mobile #"

$ws.Range("G7").Value2 = "1 2185654019"
$gf = $ws.Range("G7").Font
$gf.Name = "Segoe UI"
$gf.Bold = $true
$gf.Color = 8353126

# Row 7 is visually taller (bigger content row for bot #2)
$ws.Rows("7:7").RowHeight = 16.8

# Clear clipboard marching-ants marker
$excel.CutCopyMode = 0

# Restore the selection to where the author ended up
$ws.Range("H9").Select()
